$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve column D as text so numeric-looking price strings (e.g. "614.97")
# are not silently converted to Excel numbers when assigned below.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '63.360.49'
$ws.Range('D3').Value = '2.677.21'
$ws.Range('E3').Value = '  +4.10%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '614.97'
$ws.Range('E5').Value = '  +5.30%  '
$ws.Range('D6').Value = '143.79'
$ws.Range('E6').Value = '  -0.14%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('D9').Value = '2.676.11'
$ws.Range('E9').Value = '  +4.12%  '
$ws.Range('E10').Value = '  +0.91%  '
$ws.Range('E11').Value = '  +0.77%  '
$ws.Range('E12').Value = '  +0.59%  '
$ws.Range('D13').Value = '0.362'
$ws.Range('E13').Value = '  +3.72%  '
$ws.Range('D14').Value = '27.43'
$ws.Range('E14').Value = '  +1.73%  '
$ws.Range('D15').Value = '3.157.56'
$ws.Range('E15').Value = '  +4.12%  '
$ws.Range('D16').Value = '63.237.26'
$ws.Range('E16').Value = '  +0.60%  '
$ws.Range('E17').Value = '  +0.33%  '
$ws.Range('D18').Value = '2.688.83'
$ws.Range('E18').Value = '  +4.53%  '
$ws.Range('E19').Value = '  +3.87%  '
$ws.Range('D20').Value = '342.68'
$ws.Range('E20').Value = '  +0.70%  '
$ws.Range('D21').Value = '4.42'
$ws.Range('E21').Value = '  +2.14%  '
$ws.Range('D22').Value = '6.87'
$ws.Range('E22').Value = '  +3.77%  '
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('E24').Value = '  -0.52%  '
$ws.Range('E25').Value = '  +3.56%  '
$ws.Range('E26').Value = '  -3.72%  '
$ws.Range('D27').Value = '8.68'
$ws.Range('E27').Value = '  +5.40%  '
$ws.Range('E28').Value = '  -0.26%  '
$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D29').Value = '540.48'
$ws.Range('E29').Value = '  +16.93%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.11%  '
$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').Value = '7.92'
$ws.Range('E31').Value = '  -0.66%  '
$ws.Range('D32').Value = '2.08'
$ws.Range('E32').Value = '  +7.61%  '
$ws.Range('E33').Value = '  +8.86%  '
$ws.Range('E34').Value = '  +1.48%  '
$ws.Range('D35').Value = '172.14'
$ws.Range('E35').Value = '  -2.53%  '
$ws.Range('D36').Value = '5.18'
$ws.Range('E36').Value = '  +14.14%  '
$ws.Range('E37').Value = '  +1.78%  '
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  -0.13%  '
$ws.Range('D39').Value = '19.24'
$ws.Range('E39').Value = '  +2.23%  '
$ws.Range('E40').Value = '  +11.01%  '
$ws.Range('D41').Value = '177.14'
$ws.Range('E41').Value = '  +12.16%  '
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('E43').Value = '  +2.15%  '
$ws.Range('D44').Value = '22.31'
$ws.Range('E44').Value = '  +5.12%  '
$ws.Range('D45').Value = '0.0571'
$ws.Range('E45').Value = '  +7.06%  '
$ws.Range('D46').Value = '0.637'
$ws.Range('E46').Value = '  +0.90%  '
$ws.Range('D47').Value = '0.0964'
$ws.Range('E47').Value = '  +0.57%  '
$ws.Range('E48').Value = '  +2.27%  '
$ws.Range('D49').Value = '18.85'
$ws.Range('E49').Value = '  +4.83%  '
$ws.Range('E50').Value = '  +4.64%  '
$ws.Range('D51').Value = '11.30'
$ws.Range('E51').Value = '  -0.87%  '
